$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.710.17'
$ws.Range("E2").Value = '  -0.85%  '

# Row 3
$ws.Range("D3").Value = '2.294.40'
$ws.Range("E3").Value = '  +0.47%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.13%  '

# Row 7
$ws.Range("E7").Value = '  -0.63%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.10%  '

# Row 11
$ws.Range("E11").Value = '  -0.32%  '

# Row 12
$ws.Range("E12").Value = '  -3.60%  '

# Row 13
$ws.Range("E13").Value = '  -0.55%  '

# Row 14
$ws.Range("D14").Value = '2.636.71'
$ws.Range("E14").Value = '  +0.40%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.851'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.77%  '

# Row 17
$ws.Range("D17").Value = '2.296.88'
$ws.Range("E17").Value = '  +1.55%  '

# Row 18
$ws.Range("D18").Value = '43.582.12'
$ws.Range("E18").Value = '  -1.29%  '

# Row 19
$ws.Range("E19").Value = '  +1.83%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.75%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '

# Row 22
$ws.Range("E22").Value = '  +4.86%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -10.15%  '

# Row 25
$ws.Range("E25").Value = '  +0.01%  '

# Row 26
$ws.Range("E26").Value = '  -0.63%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.39%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.31%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.55%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0886'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.83%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.126'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("E36").Value = '  -4.93%  '

# Row 37
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.36%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.71%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.235'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.98%  '

# Row 41
$ws.Range("E41").Value = '  +4.15%  '

# Row 42
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.75%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.34%  '

# Row 44
$ws.Range("E44").Value = '  +2.49%  '

# Row 45
$ws.Range("B45").Value = 'THORChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.77%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.06%  '

# Row 47
$ws.Range("E47").Value = '  -0.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.80%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '

# Row 50
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.431'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.72%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.515.56'
$ws.Range("E51").Value = '  +0.18%  '
